# Statement update: append new transaction rows (76-87), keep style of
# blank E100:E102 consistent with neighbours, and trim three now-unused
# trailing rows (368-370) from the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New transaction rows 76-87 ------------------------------------------
$rows = @(
    @{ Row = 76; Name = "Tsepo";     Amount = 3050; Actual = "2025-08-25"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 77; Name = "Mhlengi";   Amount = 1050; Actual = "2025-09-01"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 78; Name = "Ntoko";     Amount = 1050; Actual = "2025-08-28"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 79; Name = "Bhodloza";  Amount = 0;    Actual = "2025-06-27"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 80; Name = "Msizi";     Amount = 0;    Actual = "2025-07-05"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 81; Name = "Piwe";      Amount = 550;  Actual = "2025-08-29"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 82; Name = "Mshagmor";  Amount = 1050; Actual = "2025-09-05"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 83; Name = "Thabo";     Amount = 1050; Actual = "2025-08-29"; MonthEnd = "2025-08-31"; Type = "Contribution" },
    @{ Row = 84; Name = "Thabo";     Amount = 300;  Actual = "2025-08-29"; MonthEnd = "2025-08-31"; Type = "Isipheko" },
    @{ Row = 85; Name = "Thabo";     Amount = 50;   Actual = "2025-08-29"; MonthEnd = "2025-07-31"; Type = "Contribution" },
    @{ Row = 86; Name = "Lunga";     Amount = 1050; Actual = "2025-08-20"; MonthEnd = "2025-11-30"; Type = "Contribution" },
    @{ Row = 87; Name = "Lunga";     Amount = 1050; Actual = "2025-09-02"; MonthEnd = "2025-12-31"; Type = "Contribution" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Amount
    $ws.Cells.Item($r.Row, 3).Value = $r.Actual
    $ws.Cells.Item($r.Row, 4).Value = $r.MonthEnd
    $ws.Cells.Item($r.Row, 5).Value = $r.Type
}

# --- Give blank E100:E102 the same (date) style as neighbouring C/D cells
$ws.Range("C100:C102").Copy() | Out-Null
$ws.Range("E100:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Trim the now-unused trailing rows 368-370 ---------------------------
$ws.Rows("368:370").Delete() | Out-Null

# --- Restore the view: frozen header row, scrolled down a bit, selection
$ws.Range("G66").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.FreezePanes = $true
